# Fruta / hortaliza, semanal
# Insert a new data row right after the header-adjacent block, at row 103,
# shifting all the existing rows 103..146 down to 104..147, and populate the
# newly inserted row with the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103 - this shifts rows 103:146 down to 104:147
# and keeps everything else (headers, earlier rows) untouched.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new record. It is the same
# kind of record as its (new) neighbour at row 104 (same market/product/etc.),
# just a different week, quality-"Primera", and updated price/volume figures.
$ws.Cells.Item(103, 1).Value = 11
$ws.Cells.Item(103, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(103, 3).Value = "Bíobío"
$ws.Cells.Item(103, 4).Value2 = 45006
$ws.Cells.Item(103, 4).NumberFormat = $ws.Cells.Item(104, 4).NumberFormat
$ws.Cells.Item(103, 5).Value = 8
$ws.Cells.Item(103, 6).Value = "Fruta"
$ws.Cells.Item(103, 7).Value = 100101
$ws.Cells.Item(103, 8).Value = "Berries"
$ws.Cells.Item(103, 9).Value = 100101001
$ws.Cells.Item(103, 10).Value = "Arándano (blue)"
$ws.Cells.Item(103, 11).Value = "Sin especificar"
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 220
$ws.Cells.Item(103, 14).Value = 3500
$ws.Cells.Item(103, 15).Value = 4000
$ws.Cells.Item(103, 16).Value = 3773
$ws.Cells.Item(103, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(103, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(103, 19).Value = 1886
$ws.Cells.Item(103, 20).Value = 2
